# Update 2p0. Convention change to support multi-axle vehicles
#
# - Insert a new "Truck_Amandla" sheet (copy of Trailer_Elula's template)
#   right before "Trailer_Elula".
# - Append a new "Trailer_Kumanzi" sheet (copy of Trailer_Thwala's template)
#   at the end of the workbook, and make it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- Truck_Amandla: new sheet inserted before Trailer_Elula ---------------
$elula = $wb.Worksheets.Item("Trailer_Elula")
[void]$elula.Copy($elula)
$truck = $wb.Worksheets.Item("Trailer_Elula (2)")
$truck.Name = "Truck_Amandla"

$truck.Range("H3").Value = "Truck_Amandla"
$truck.Range("H6").Value = 0.43
$truck.Range("F9").Value = -1.2
$truck.Range("G9").Value = 0
$truck.Range("H9").Value = 1.1

$truck.Activate()
[void]$truck.Range("H5:H9").Select()

# --- Trailer_Kumanzi: new sheet appended after Trailer_Thwala --------------
$thwala = $wb.Worksheets.Item("Trailer_Thwala")
[void]$thwala.Copy($null, $thwala)
$kumanzi = $wb.Worksheets.Item("Trailer_Thwala (2)")
$kumanzi.Name = "Trailer_Kumanzi"

$kumanzi.Range("H3").Value = "Trailer_Kumanzi"
$kumanzi.Range("H6").Value = 0.43
$kumanzi.Range("F9").Value = 5
$kumanzi.Range("G9").Value = 0
$kumanzi.Range("H9").Value = 2

# Trailer_Kumanzi ends up as the active/selected sheet.
$kumanzi.Activate()
[void]$kumanzi.Range("J20").Select()
